$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "LP1912" (main schedule sheet)
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("LP1912")

# Header text updates
$ws1.Range("A2").Value = "Última actualización: 14:45:56"
$ws1.Range("A3").Value = "Total filas: 213"

# --- Swap row 58 <-> row 59 -------------------------------------------------
$a58 = $ws1.Cells.Item(58,1).Value()
$b58 = $ws1.Cells.Item(58,2).Value()
$c58 = $ws1.Cells.Item(58,3).Value()
$d58 = $ws1.Cells.Item(58,4).Value()
$e58 = $ws1.Cells.Item(58,5).Value()

$a59 = $ws1.Cells.Item(59,1).Value()
$b59 = $ws1.Cells.Item(59,2).Value()
$c59 = $ws1.Cells.Item(59,3).Value()
$d59 = $ws1.Cells.Item(59,4).Value()
$e59 = $ws1.Cells.Item(59,5).Value()

$ws1.Cells.Item(58,1).Value = $a59
$ws1.Cells.Item(58,2).Value = $b59
$ws1.Cells.Item(58,3).Value = $c59
$ws1.Cells.Item(58,4).Value = $d59
$ws1.Cells.Item(58,5).Value = $e59

$ws1.Cells.Item(59,1).Value = $a58
$ws1.Cells.Item(59,2).Value = $b58
$ws1.Cells.Item(59,3).Value = $c58
$ws1.Cells.Item(59,4).Value = $d58
$ws1.Cells.Item(59,5).Value = $e58

# --- Swap row 111 <-> row 113 (row 112 untouched) --------------------------
$a111 = $ws1.Cells.Item(111,1).Value()
$b111 = $ws1.Cells.Item(111,2).Value()
$c111 = $ws1.Cells.Item(111,3).Value()
$d111 = $ws1.Cells.Item(111,4).Value()
$e111 = $ws1.Cells.Item(111,5).Value()

$a113 = $ws1.Cells.Item(113,1).Value()
$b113 = $ws1.Cells.Item(113,2).Value()
$c113 = $ws1.Cells.Item(113,3).Value()
$d113 = $ws1.Cells.Item(113,4).Value()
$e113 = $ws1.Cells.Item(113,5).Value()

$ws1.Cells.Item(111,1).Value = $a113
$ws1.Cells.Item(111,2).Value = $b113
$ws1.Cells.Item(111,3).Value = $c113
$ws1.Cells.Item(111,4).Value = $d113
$ws1.Cells.Item(111,5).Value = $e113

$ws1.Cells.Item(113,1).Value = $a111
$ws1.Cells.Item(113,2).Value = $b111
$ws1.Cells.Item(113,3).Value = $c111
$ws1.Cells.Item(113,4).Value = $d111
$ws1.Cells.Item(113,5).Value = $e111

# --- Swap row 137 <-> row 138 ------------------------------------------------
$a137 = $ws1.Cells.Item(137,1).Value()
$b137 = $ws1.Cells.Item(137,2).Value()
$c137 = $ws1.Cells.Item(137,3).Value()
$d137 = $ws1.Cells.Item(137,4).Value()
$e137 = $ws1.Cells.Item(137,5).Value()

$a138 = $ws1.Cells.Item(138,1).Value()
$b138 = $ws1.Cells.Item(138,2).Value()
$c138 = $ws1.Cells.Item(138,3).Value()
$d138 = $ws1.Cells.Item(138,4).Value()
$e138 = $ws1.Cells.Item(138,5).Value()

$ws1.Cells.Item(137,1).Value = $a138
$ws1.Cells.Item(137,2).Value = $b138
$ws1.Cells.Item(137,3).Value = $c138
$ws1.Cells.Item(137,4).Value = $d138
$ws1.Cells.Item(137,5).Value = $e138

$ws1.Cells.Item(138,1).Value = $a137
$ws1.Cells.Item(138,2).Value = $b137
$ws1.Cells.Item(138,3).Value = $c137
$ws1.Cells.Item(138,4).Value = $d137
$ws1.Cells.Item(138,5).Value = $e137

# --- Insert two new rows before the old row 207 (pushes 207-214 to 209-216) -
$ws1.Rows("207:208").Insert()

$ws1.Cells.Item(207,1).Value = "14:45:56"
$ws1.Cells.Item(207,2).Value = "15:50"
$ws1.Cells.Item(207,3).Value = "27_EL RETIRO"
$ws1.Cells.Item(207,4).Value = 65
$ws1.Cells.Item(207,5).Value = "LP1912"

$ws1.Cells.Item(208,1).Value = "14:45:56"
$ws1.Cells.Item(208,2).Value = "15:52"
$ws1.Cells.Item(208,3).Value = "10_OLMOS"
$ws1.Cells.Item(208,4).Value = 67
$ws1.Cells.Item(208,5).Value = "LP1912"

# --- Append two new rows at the end (217 and 218) ---------------------------
$ws1.Cells.Item(217,1).Value = "14:45:56"
$ws1.Cells.Item(217,2).Value = "16:33"
$ws1.Cells.Item(217,3).Value = "83_ALUAR"
$ws1.Cells.Item(217,4).Value = 108
$ws1.Cells.Item(217,5).Value = "LP1912"

$ws1.Cells.Item(218,1).Value = "14:45:56"
$ws1.Cells.Item(218,2).Value = "16:40"
$ws1.Cells.Item(218,3).Value = "225_GOMEZ"
$ws1.Cells.Item(218,4).Value = 115
$ws1.Cells.Item(218,5).Value = "LP1912"

# ---------------------------------------------------------------------------
# Sheet "LP1912-215"
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("LP1912-215")
$ws2.Range("A2").Value = "Última actualización: 14:45:56"

# ---------------------------------------------------------------------------
# Sheet "6203-6173"
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("6203-6173")
$ws3.Range("A2").Value = "Última actualización: 14:45:56"
$ws3.Range("A3").Value = "Total filas: 32"

# Insert one new row before the old row 36 (pushes it to row 37)
$ws3.Rows("36:36").Insert()

$ws3.Cells.Item(36,1).Value = "14:45:56"
$ws3.Cells.Item(36,2).Value = "16:05"
$ws3.Cells.Item(36,3).Value = "215C_LA PLATA"
$ws3.Cells.Item(36,4).Value = 80
$ws3.Cells.Item(36,5).Value = "L6203"
